# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 8601
$ws1.Range("F3").Value  = 74
$ws1.Range("F6").Value  = 1371
$ws1.Range("F7").Value  = 130
$ws1.Range("F8").Value  = 27
$ws1.Range("F10").Value = 9348
$ws1.Range("F12").Value = 96
$ws1.Range("F13").Value = 219
$ws1.Range("F15").Value = 355
$ws1.Range("F16").Value = 6324
$ws1.Range("F17").Value = 1062
$ws1.Range("F18").Value = 81
$ws1.Range("F19").Value = 45
$ws1.Range("F20").Value = 129

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 8601
$ws4.Range("F3").Value  = 74
$ws4.Range("F6").Value  = 1371
$ws4.Range("F7").Value  = 130
$ws4.Range("F8").Value  = 27
$ws4.Range("F12").Value = 9348
$ws4.Range("F14").Value = 96
$ws4.Range("F15").Value = 219
$ws4.Range("F17").Value = 355
$ws4.Range("F18").Value = 6324
$ws4.Range("F19").Value = 1062
$ws4.Range("F20").Value = 81
$ws4.Range("F21").Value = 45
$ws4.Range("F22").Value = 129

$wb.Save()
